$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "path"
$ws.Range("C1").Value = "short"
$ws.Range("D1").Value = "definition"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "card."
$ws.Range("G1").Value = "binding"

$ws.Range("A2").Value = "EHDSCondition"
$ws.Range("B2").Value = "EHDSCondition"
$ws.Range("C2").Value = "Condition model"
$ws.Range("D2").Value = "EHDS refined base model for A clinical condition, problem, diagnosis, or other event, situation, issue, or clinical concept that has risen to a level of concern."
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "0..*"
$ws.Range("G2").Value = ""

$ws.Range("A3").Value = "EHDSCondition.identifier"
$ws.Range("B3").Value = "EHDSCondition.identifier"
$ws.Range("C3").Value = "Condition identifier"
$ws.Range("D3").Value = "Condition identifier"
$ws.Range("E3").Value = "Identifier"
$ws.Range("F3").Value = "0..*"
$ws.Range("G3").Value = ""

$ws.Range("A4").Value = "EHDSCondition.subject"
$ws.Range("B4").Value = "EHDSCondition.subject"
$ws.Range("C4").Value = "Indicates the patient or group who the condition record is associated with."
$ws.Range("D4").Value = "Indicates the patient or group who the condition record is associated with."
$ws.Range("E4").Value = "EHDSPatient"
$ws.Range("F4").Value = "1..1"
$ws.Range("G4").Value = ""

$ws.Range("A5").Value = "EHDSCondition.description"
$ws.Range("B5").Value = "EHDSCondition.description"
$ws.Range("C5").Value = "Condition specification in narrative form"
$ws.Range("D5").Value = "Condition specification in narrative form"
$ws.Range("E5").Value = "string"
$ws.Range("F5").Value = "0..1"
$ws.Range("G5").Value = ""

$ws.Range("A6").Value = "EHDSCondition.code"
$ws.Range("B6").Value = "EHDSCondition.code"
$ws.Range("C6").Value = "Code identifying the condition, problem or diagnosis"
$ws.Range("D6").Value = "Code identifying the condition, problem or diagnosis"
$ws.Range("E6").Value = "CodeableConcept"
$ws.Range("F6").Value = "0..1"
$ws.Range("G6").Value = "{'strength': 'preferred', 'description': 'ICD-10*, SNOMED CT, Orphacode if rare disease is diagnosed'}"

$ws.Range("A7").Value = "EHDSCondition.onsetDate"
$ws.Range("B7").Value = "EHDSCondition.onsetDate"
$ws.Range("C7").Value = "Onset date of a problem/condition"
$ws.Range("D7").Value = "Onset date of a problem/condition"
$ws.Range("E7").Value = "dateTime"
$ws.Range("F7").Value = "0..1"
$ws.Range("G7").Value = ""

$ws.Range("A8").Value = "EHDSCondition.endDate"
$ws.Range("B8").Value = "EHDSCondition.endDate"
$ws.Range("C8").Value = "The date or estimated date that the condition resolved or went into remission."
$ws.Range("D8").Value = "The date or estimated date that the condition resolved or went into remission."
$ws.Range("E8").Value = "dateTime"
$ws.Range("F8").Value = "0..1"
$ws.Range("G8").Value = ""

$ws.Range("A9").Value = "EHDSCondition.category"
$ws.Range("B9").Value = "EHDSCondition.category"
$ws.Range("C9").Value = "Category or categories of the problem."
$ws.Range("D9").Value = "Category or categories of the problem."
$ws.Range("E9").Value = "CodeableConcept"
$ws.Range("F9").Value = "0..*"
$ws.Range("G9").Value = ""

$ws.Range("A10").Value = "EHDSCondition.clinicalStatus"
$ws.Range("B10").Value = "EHDSCondition.clinicalStatus"
$ws.Range("C10").Value = "Status of the condition/problem (active, resolved, inactive, ...)"
$ws.Range("D10").Value = "Status of the condition/problem (active, resolved, inactive, ...)"
$ws.Range("E10").Value = "CodeableConcept"
$ws.Range("F10").Value = "0..1"
$ws.Range("G10").Value = "{'strength': 'preferred', 'description': 'HL7 Condition-clinical'}"

$ws.Range("A11").Value = "EHDSCondition.resolutionCircumstances[x]"
$ws.Range("B11").Value = "EHDSCondition.resolutionCircumstances[x]"
$ws.Range("C11").Value = "Describes the reason for which the status of the problem changed from current to inactive (e.g. surgical procedure, medical treatment, etc.)."
$ws.Range("D11").Value = "This field includes free text if the resolution circumstances are not already included in other fields such as surgical procedure, medical device, etc., e.g. hepatic cystectomy (this will be the resolution circumstances for the problem `"hepatic cyst`" and will be included in surgical procedures)."
$ws.Range("E11").Value = "CodeableReference"
$ws.Range("F11").Value = "0..*"
$ws.Range("G11").Value = ""

$ws.Range("A12").Value = "EHDSCondition.severity"
$ws.Range("B12").Value = "EHDSCondition.severity"
$ws.Range("C12").Value = "A subjective assessment of the severity of the condition as evaluated by the clinician."
$ws.Range("D12").Value = "A subjective assessment of the severity of the condition as evaluated by the clinician."
$ws.Range("E12").Value = "CodeableConcept"
$ws.Range("F12").Value = "0..1"
$ws.Range("G12").Value = "{'strength': 'preferred', 'description': 'HL7 Condition-severity'}"

$ws.Range("A13").Value = "EHDSCondition.anatomicLocation[x]"
$ws.Range("B13").Value = "EHDSCondition.anatomicLocation[x]"
$ws.Range("C13").Value = "The anatomical location including laterality where this condition manifests itself."
$ws.Range("D13").Value = "The anatomical location including laterality where this condition manifests itself."
$ws.Range("E13").Value = "CodeableConcept"
$ws.Range("F13").Value = "0..*"
$ws.Range("G13").Value = "{'strength': 'preferred', 'description': 'SNOMED CT'}"

$ws.Range("A14").Value = "EHDSCondition.stage"
$ws.Range("B14").Value = "EHDSCondition.stage"
$ws.Range("C14").Value = "Stage/grade usually assessed formally using a specific staging/grading system. Multiple assessment systems could be used."
$ws.Range("D14").Value = "Stage/grade usually assessed formally using a specific staging/grading system. Multiple assessment systems could be used."
$ws.Range("E14").Value = "CodeableConcept"
$ws.Range("F14").Value = "0..*"
$ws.Range("G14").Value = "{'strength': 'preferred', 'description': 'e.g. TNM, ICD-O-3, Bi-Rads, Li-Rads, …'}"

$ws.Range("A15").Value = "EHDSCondition.diagnosisAssertionStatus"
$ws.Range("B15").Value = "EHDSCondition.diagnosisAssertionStatus"
$ws.Range("C15").Value = "Assertion about the certainty associated with a diagnosis. Diagnostic and/or clinical evidence of condition."
$ws.Range("D15").Value = "Assertion about the certainty associated with a diagnosis. Diagnostic and/or clinical evidence of condition."
$ws.Range("E15").Value = "CodeableConcept"
$ws.Range("F15").Value = "0..1"
$ws.Range("G15").Value = "{'strength': 'preferred', 'description': 'HL7 Condition-ver-status'}"

$ws.Range("A16").Value = "EHDSCondition.asserter"
$ws.Range("B16").Value = "EHDSCondition.asserter"
$ws.Range("C16").Value = "The asserter of the condition"
$ws.Range("D16").Value = "The asserter of the condition"
$ws.Range("E16").Value = "EHDSHealthProfessional"
$ws.Range("F16").Value = "0..1"
$ws.Range("G16").Value = ""

$ws.Range("A17").Value = "EHDSCondition.assertedDate"
$ws.Range("B17").Value = "EHDSCondition.assertedDate"
$ws.Range("C17").Value = "Date and time of the diagnosis assertion"
$ws.Range("D17").Value = "Date and time of the diagnosis assertion"
$ws.Range("E17").Value = "dateTime"
$ws.Range("F17").Value = "0..1"
$ws.Range("G17").Value = ""

$ws.Range("A18:G19").EntireRow.Delete()